$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update column C (Accident Id) for rows 2-6 ---
$ws.Range("C2").Value = "A-2827637"
$ws.Range("C3").Value = "A-2827637"
$ws.Range("C4").Value = "A-2827637"
$ws.Range("C5").Value = "A-2827637"
$ws.Range("C6").Value = "A-2827637"

# --- Update column D (Image_link) filenames ---
$ws.Range("D2").Value = "00001.png"
$ws.Range("D3").Value = "0002.png"
$ws.Range("D4").Value = "0003.png"
$ws.Range("D5").Value = "0004.png"
$ws.Range("D6").Value = "0005.png"

# --- Add two new rows of data (7 and 8) ---
$ws.Range("A7").Value = 6
$ws.Range("B7").Value = 2
$ws.Range("C7").Value = "A-2827637"
$ws.Range("D7").Value = "0006.png"
$ws.Range("E7").Value = "no class"

$ws.Range("A8").Value = 7
$ws.Range("B8").Value = 2
$ws.Range("C8").Value = "A-2827637"
$ws.Range("D8").Value = "0007.png"
$ws.Range("E8").Value = 14

# --- Remove the old stray row 16 (E16 = "A") ---
$ws.Range("A16:E16").ClearContents()

# --- Update the selection to match the authored state ---
$ws.Range("C11").Select()
